$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "panel" header cell (E1) onto the
# new "time_taken" header cell (F1) so it reuses the same bold/bordered style.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

$times = @(
    "2021-10-05 13:40:36.685498",
    "2021-10-05 13:40:36.685511",
    "2021-10-05 13:40:36.685515",
    "2021-10-05 13:40:36.685518",
    "2021-10-05 13:40:36.685522",
    "2021-10-05 13:40:36.685525",
    "2021-10-05 13:40:36.685528",
    "2021-10-05 13:40:36.685531",
    "2021-10-05 13:40:36.685534",
    "2021-10-05 13:40:36.685538",
    "2021-10-05 13:40:36.685541",
    "2021-10-05 13:40:36.685544",
    "2021-10-05 13:40:36.685547",
    "2021-10-05 13:40:36.685550",
    "2021-10-05 13:40:36.685553",
    "2021-10-05 13:40:36.685556",
    "2021-10-05 13:40:36.685560"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
